# Updated cryptos list on Tue Nov 19 16:53:07 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and fixes the ordering/data for two row pairs whose rank swapped
# (Avalanche/WrappedBTC at rows 14-15, VeChain/ImmutableX at rows 49-50).
#
# Several Price values are plain decimal-looking text (e.g. "241.40",
# "0.0000256") that must stay stored as literal text, not be coerced into a
# floating point number (which would silently drop trailing zeros / change
# precision). Prefixing with a leading apostrophe forces Excel to treat the
# assignment as text, matching how these inline-string cells already exist
# in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "92.368.26"
$ws.Range("E2").Value = "  +1.07%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.110.29"
$ws.Range("E3").Value = "  -1.68%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.42%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'241.40"
$ws.Range("E5").Value = "  +0.59%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'614.27"
$ws.Range("E6").Value = "  -1.15%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -3.12%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "'0.401"
$ws.Range("E8").Value = "  +8.41%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.04%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "3.106.05"
$ws.Range("E10").Value = "  -1.78%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.731"
$ws.Range("E11").Value = "  -1.70%  "

# Row 12 - Toncoin-like / next coin (only Volume changed)
$ws.Range("E12").Value = "  -0.88%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  +2.98%  "

# Row 14 - now WrappedBTC (was Avalanche)
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "92.158.97"
$ws.Range("E14").Value = "  +0.85%  "

# Row 15 - now Avalanche (was WrappedBTC)
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'34.42"
$ws.Range("E15").Value = "  -2.92%  "

# Row 16 - Toncoin
$ws.Range("D16").Value = "'5.51"
$ws.Range("E16").Value = "  +0.18%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.689.07"
$ws.Range("E17").Value = "  -1.57%  "

# Row 18 - next coin
$ws.Range("D18").Value = "3.070.54"
$ws.Range("E18").Value = "  -2.97%  "

# Row 19 - (only Volume changed)
$ws.Range("E19").Value = "  -1.42%  "

# Row 20
$ws.Range("D20").Value = "'14.76"
$ws.Range("E20").Value = "  -5.04%  "

# Row 21
$ws.Range("D21").Value = "'5.81"
$ws.Range("E21").Value = "  -0.44%  "

# Row 22
$ws.Range("D22").Value = "'9.40"
$ws.Range("E22").Value = "  +2.16%  "

# Row 23
$ws.Range("D23").Value = "'447.90"
$ws.Range("E23").Value = "  +0.89%  "

# Row 24
$ws.Range("D24").Value = "'0.0000202"
$ws.Range("E24").Value = "  -3.60%  "

# Row 25
$ws.Range("D25").Value = "'5.79"
$ws.Range("E25").Value = "  -0.39%  "

# Row 26
$ws.Range("D26").Value = "'87.24"
$ws.Range("E26").Value = "  +5.31%  "

# Row 27
$ws.Range("D27").Value = "'11.77"
$ws.Range("E27").Value = "  -2.38%  "

# Row 29 - (only Volume changed)
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("D30").Value = "'0.137"
$ws.Range("E30").Value = "  -0.39%  "

# Row 31
$ws.Range("D31").Value = "'0.233"
$ws.Range("E31").Value = "  -1.68%  "

# Row 32 - (only Volume changed)
$ws.Range("E32").Value = "  -1.90%  "

# Row 33
$ws.Range("D33").Value = "'9.31"
$ws.Range("E33").Value = "  +0.02%  "

# Row 34
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +25.18%  "

# Row 35
$ws.Range("D35").Value = "'8.08"
$ws.Range("E35").Value = "  +2.66%  "

# Row 36
$ws.Range("D36").Value = "'0.166"
$ws.Range("E36").Value = "  -0.06%  "

# Row 37
$ws.Range("D37").Value = "'4.28"
$ws.Range("E37").Value = "  +3.24%  "

# Row 38
$ws.Range("D38").Value = "'26.18"
$ws.Range("E38").Value = "  -1.46%  "

# Row 39 - (only Volume changed)
$ws.Range("E39").Value = "  -0.08%  "

# Row 40 - (only Volume changed)
$ws.Range("E40").Value = "  -2.15%  "

# Row 41
$ws.Range("D41").Value = "'481.40"
$ws.Range("E41").Value = "  -6.13%  "

# Row 42
$ws.Range("D42").Value = "'0.442"
$ws.Range("E42").Value = "  -5.14%  "

# Row 43
$ws.Range("D43").Value = "'3.50"
$ws.Range("E43").Value = "  +0.62%  "

# Row 44 - (only Volume changed)
$ws.Range("E44").Value = "  +3.99%  "

# Row 45 - (only Volume changed)
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("D46").Value = "'159.48"
$ws.Range("E46").Value = "  +2.33%  "

# Row 47 - (only Volume changed)
$ws.Range("E47").Value = "  -0.17%  "

# Row 48
$ws.Range("D48").Value = "'0.692"
$ws.Range("E48").Value = "  -3.46%  "

# Row 49 - now ImmutableX (was VeChain)
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "'1.37"
$ws.Range("E49").Value = "  -1.35%  "

# Row 50 - now VeChain (was ImmutableX)
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0335"
$ws.Range("E50").Value = "  +2.71%  "

# Row 51 - OKB
$ws.Range("D51").Value = "'44.13"
$ws.Range("E51").Value = "  +0.04%  "
